$d = $word.ActiveDocument

$d.Content.Find.Execute("neighboring cities: A, B and C decided", $true, $false, $false, $false, $false, $true, 1, $false, "ከተማ መስተዳድሮች:- ሀ፣ ለ፣ ሐ ወሰኑ", 2) | Out-Null
$d.Content.Find.Execute("to build an airport dividing the costs of", $true, $false, $false, $false, $false, $true, 1, $false, "የአየር ማረፊያ ለመገንባት የትግበራ ወጭውን", 2) | Out-Null
$d.Content.Find.Execute("implementation. The condition on the", $true, $false, $false, $false, $false, $true, 1, $false, "በመካፈል፡፡  መስፈርቱ", 2) | Out-Null
$d.Content.Find.Execute("choice of the most suitable place is", $true, $false, $false, $false, $false, $true, 1, $false, "በጣም ምቹ ቦታ ለመምረጥ", 2) | Out-Null
$d.Content.Find.Execute("that the sum of the distances from each", $true, $false, $false, $false, $false, $true, 1, $false, "ከሶስቱ ከተሞች እስከ አየር ማረፊያው ያሉት እርቀቶች", 2) | Out-Null
$d.Content.Find.Execute("city to the airport is as small as", $true, $false, $false, $false, $false, $true, 1, $false, "ተደምረው አነስተኛ እንድሆን", 2) | Out-Null
$d.Content.Find.Execute("possible. The team of experts in charge", $true, $false, $false, $false, $false, $true, 1, $false, "በተቻለ መጠን፡፡ ሃላፊነቱን የወሰደው የባለሙያዎች ቡድን", 2) | Out-Null
$d.Content.Find.Execute("of the work has created a model to get", $true, $false, $false, $false, $false, $true, 1, $false, "ሞደል ፈጠረ", 2) | Out-Null
$d.Content.Find.Execute("a preliminary idea of where to place the", $true, $false, $false, $false, $false, $true, 1, $false, "መነሻ ሃሳብ ለማገኘት ቦታው", 2) | Out-Null
$d.Content.Find.Execute("structure. At their disposal there are", $true, $false, $false, $false, $false, $true, 1, $false, "የት መሆን እንዳለበት፡፡ ማስወገጃቸው ጋር", 2) | Out-Null
$d.Content.Find.Execute("some snails a big metal ring and a long", $true, $false, $false, $false, $false, $true, 1, $false, "ትንንሽ ሚሰማሮች፣ትልቅ የብረት ቀለበት እና  ረጅም", 2) | Out-Null
$d.Content.Find.Execute("string.", $true, $false, $false, $false, $false, $true, 1, $false, "ተለጣጭ ገመድ አሉ፡፡", 2) | Out-Null
$d.Content.Find.Execute("Explain how the team can manage to use", $true, $false, $false, $false, $false, $true, 1, $false, "አብራሩ ቡድኑ እንደት ማደረግ እንዳለበት", 2) | Out-Null
$d.Content.Find.Execute("the materials to tell approximately the", $true, $false, $false, $false, $false, $true, 1, $false, "መሳሪያዎቹን ተቀራራቢ", 2) | Out-Null
$d.Content.Find.Execute("ideal location of the airport. Imagine", $true, $false, $false, $false, $false, $true, 1, $false, "ሃሳባዊ የአየር መንገዱን ለመናገር፡፡ አስቡ", 2) | Out-Null
$d.Content.Find.Execute("that the cities are placed at the", $true, $false, $false, $false, $false, $true, 1, $false, "ከተሞቹ", 2) | Out-Null
$d.Content.Find.Execute("vertices of a triangle which is", $true, $false, $false, $false, $false, $true, 1, $false, "የሶስት መአዘን ጠርዞች ላይ እንዳሉ፣ ይህም", 2) | Out-Null
$d.Content.Find.Execute("obviously reproduced in scale as", $true, $false, $false, $false, $false, $true, 1, $false, "በግልጽ በልኬቱ እንደተሰራው", 2) | Out-Null
$d.Content.Find.Execute("shown in figure. This is one possible", $true, $false, $false, $false, $false, $true, 1, $false, "ምስሉ ላይ እነደሚታየው ነው፡፡ ይኸ አንዱ  አማራጭ", 2) | Out-Null
$d.Content.Find.Execute("setting the rope starts from one nail,", $true, $false, $false, $false, $false, $true, 1, $false, "መንገድ ገመዱ ከአንዱ ሚስማር ጀምሮ", 2) | Out-Null
$d.Content.Find.Execute("goes inside the ring, goes around the", $true, $false, $false, $false, $false, $true, 1, $false, "በቀለበቱ ውስጥ አልፎ፣ በሌለኛው", 2) | Out-Null
$d.Content.Find.Execute("other nail, the third nail, inside the", $true, $false, $false, $false, $false, $true, 1, $false, "ሚስማር ዞሮ፣ በሶስተኛው ሚስማር፣", 2) | Out-Null
$d.Content.Find.Execute("ring again and now you can just pull the", $true, $false, $false, $false, $false, $true, 1, $false, "እንደገና በቀለበቱ ውስጥ አሁን", 2) | Out-Null
$d.Content.Find.Execute("rope in order to find the point that", $true, $false, $false, $false, $false, $true, 1, $false, "ገመዱን መሳብ ትችላላችሁ", 2) | Out-Null
$d.Content.Find.Execute("you're looking for. In order to reach the", $true, $false, $false, $false, $false, $true, 1, $false, "የምተፈልጉትን ነጥብ ለማገኘት፡፡ ነጥቡ ላይ ለመድረስ፤", 2) | Out-Null
$d.Content.Find.Execute("point, we have to move the rope a bit", $true, $false, $false, $false, $false, $true, 1, $false, "ገመዱን በትንሹ ማንቀሳቀስ አለብን", 2) | Out-Null
$d.Content.Find.Execute("by the materials that we are using but", $true, $false, $false, $false, $false, $true, 1, $false, "ከተጠቀምንባቸው መሳሪያዎች ስለሚኖር ግን", 2) | Out-Null
$d.Content.Find.Execute("after a while you'll reach a position from", $true, $false, $false, $false, $false, $true, 1, $false, "ከተወሰነ ጊዜ በኋላ ገመዱ መንቀሳቀስ የማይችልበት", 2) | Out-Null
$d.Content.Find.Execute("which the ring doesn't move anymore,", $true, $false, $false, $false, $false, $true, 1, $false, "ቦታ ትደርሳላችሁ፣", 2) | Out-Null
$d.Content.Find.Execute("which is more or less this one. And as", $true, $false, $false, $false, $false, $true, 1, $false, "ከሞላ ገደል ይህ ቦታ ነው፡፡ እና", 2) | Out-Null
$d.Content.Find.Execute("you can see the three distances", $true, $false, $false, $false, $false, $true, 1, $false, "እንደምታዩት የሶስቱ እርቀቶች", 2) | Out-Null
$d.Content.Find.Execute("between the ring and the nails are", $true, $false, $false, $false, $false, $true, 1, $false, "ከመወጠሪያውና ከሚስማሮቹ", 2) | Out-Null
$d.Content.Find.Execute("placed more or less 120 degrees from one", $true, $false, $false, $false, $false, $true, 1, $false, "120 ዲገሪ ነወ ከሞላ ጎደል ከአንዱ እስከ አንዱ፣ ", 2) | Out-Null
$d.Content.Find.Execute("another which is 1/3 of a circumference,", $true, $false, $false, $false, $false, $true, 1, $false, "ይህም የዙሪያው 1/3 ኛ ነው፣", 2) | Out-Null
$d.Content.Find.Execute("and that's the point that we're looking", $true, $false, $false, $false, $false, $true, 1, $false, "እና ይህ የምንፈልገው ነጥብ ነው፡-", 2) | Out-Null
$d.Content.Find.Execute("for: the minimum distance between the", $true, $false, $false, $false, $false, $true, 1, $false, "ትንሹ እርቀት በ", 2) | Out-Null
$d.Content.Find.Execute("nails and the airport when you sum it", $true, $false, $false, $false, $false, $true, 1, $false, "ሚስማሮቹ እና በአየር መንገዱ መካከል ስትደምሩት", 2) | Out-Null
$d.Content.Find.Execute("[Music]", $true, $false, $false, $false, $false, $true, 1, $false, "[ሙዚቃ]", 2) | Out-Null
